$d = $word.ActiveDocument

$d.Content.Find.Execute("August 06, 2020", $true, $false, $false, $false, $false,
                         $true, 1, $false, "August 18, 2020", 2)

$d.Content.Find.Execute("Davao Sugar Central Company, Inc.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Cotabato Sugar Central Company, Inc.", 2)

$d.Content.Find.Execute("5/F, Filinvest Bldg., No. 79 EDSA Highway Hills, Mandaluyong City", $true, $false, $false, $false, $false,
                         $true, 1, $false, "6/F Filinvest Bldg., No. 79 EDSA, Highway Hills, Mandaluyong City", 2)

$d.Content.Find.Execute("Dear Mr. Gotianun:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Dear Pres. Gotianun:", 2)
